$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 43 / 44 swap (VeChain <-> ApeXProtocol) plus new values ---
$ws.Range("B43").Value = "VeChain"
$ws.Range("C43").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.0290"
$ws.Range("E43").Value = "  +4.71%  "

$ws.Range("B44").Value = "ApeXProtocol"
$ws.Range("C44").Value = "https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
$ws.Range("D44").Value = "2.26"
$ws.Range("E44").Value = "  -1.27%  "

# --- Row 51 (HuobiToken -> THORChain) ---
$ws.Range("B51").Value = "THORChain"
$ws.Range("C51").Value = "https://coinranking.com/coin/ybmU-kKU+thorchain-rune"
$ws.Range("D51").Value = "4.57"
$ws.Range("E51").Value = "  +2.21%  "

# --- Price (D) / Volume(1h) (E) updates for remaining rows ---
$ws.Range("D2").Value = "43.179.84"
$ws.Range("E2").Value = "  +2.36%  "

$ws.Range("D3").Value = "2.316.08"
$ws.Range("E3").Value = "  +1.88%  "

$ws.Range("E4").Value = "  +0.04%  "

$ws.Range("D5").Value = "302.42"
$ws.Range("E5").Value = "  +1.10%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "101.70"
$ws.Range("E6").Value = "  +6.80%  "

$ws.Range("D7").Value = "0.504"

$ws.Range("E8").Value = "  +0.01%  "

$ws.Range("D9").Value = "0.517"
$ws.Range("E9").Value = "  +5.55%  "

$ws.Range("D10").Value = "36.13"
$ws.Range("E10").Value = "  +9.19%  "

$ws.Range("D11").Value = "0.0795"
$ws.Range("E11").Value = "  +1.00%  "

$ws.Range("E12").Value = "  +3.46%  "

$ws.Range("D13").Value = "17.95"
$ws.Range("E13").Value = "  +13.03%  "

$ws.Range("D14").Value = "6.91"
$ws.Range("E14").Value = "  +3.75%  "

$ws.Range("D15").Value = "2.675.33"
$ws.Range("E15").Value = "  +1.88%  "

$ws.Range("D16").Value = "2.319.30"
$ws.Range("E16").Value = "  +2.50%  "

$ws.Range("D17").Value = "0.809"
$ws.Range("E17").Value = "  +3.31%  "

$ws.Range("D18").Value = "43.077.73"
$ws.Range("E18").Value = "  +2.20%  "

$ws.Range("D19").Value = "12.63"
$ws.Range("E19").Value = "  +7.83%  "

$ws.Range("D20").Value = "6.21"
$ws.Range("E20").Value = "  +4.00%  "

$ws.Range("D22").Value = "67.87"
$ws.Range("E22").Value = "  +2.53%  "

$ws.Range("D23").Value = "236.72"
$ws.Range("E23").Value = "  +0.69%  "

$ws.Range("D24").Value = "2.21"
$ws.Range("E24").Value = "  +13.06%  "

$ws.Range("D25").Value = "2.46"
$ws.Range("E25").Value = "  +0.45%  "

$ws.Range("E26").Value = "  -0.10%  "

$ws.Range("D27").Value = "24.77"
$ws.Range("E27").Value = "  +4.25%  "

$ws.Range("D28").Value = "2.35"
$ws.Range("E28").Value = "  +5.81%  "

$ws.Range("D29").Value = "34.82"
$ws.Range("E29").Value = "  +3.11%  "

$ws.Range("D30").Value = "168.91"
$ws.Range("E30").Value = "  +0.35%  "

$ws.Range("D31").Value = "9.22"
$ws.Range("E31").Value = "  +0.60%  "

$ws.Range("E32").Value = "  +0.04%  "

$ws.Range("D33").Value = "4.75"
$ws.Range("E33").Value = "  +2.40%  "

$ws.Range("D34").Value = "5.03"
$ws.Range("E34").Value = "  +2.74%  "

$ws.Range("D35").Value = "17.35"
$ws.Range("E35").Value = "  +3.58%  "

$ws.Range("E36").Value = "  +2.55%  "

$ws.Range("E37").Value = "  +1.08%  "

$ws.Range("E38").Value = "  +4.51%  "

$ws.Range("E39").Value = "  +2.16%  "

$ws.Range("E40").Value = "  +4.41%  "

$ws.Range("E41").Value = "  +1.74%  "

$ws.Range("D42").Value = "1.989.15"

$ws.Range("D45").Value = "10.24"
$ws.Range("E45").Value = "  +6.76%  "

$ws.Range("E46").Value = "  +5.23%  "

$ws.Range("D47").Value = "17.59"
$ws.Range("E47").Value = "  -0.44%  "

$ws.Range("D48").Value = "56.24"
$ws.Range("E48").Value = "  +7.74%  "

$ws.Range("D49").Value = "2.546.79"
$ws.Range("E49").Value = "  +1.98%  "

$ws.Range("E50").Value = "  +3.40%  "
